# Update countries & provincias Spain
#
# The source workbook is a COVID-19 daily tracker ("Pais" sheet). Each
# refresh re-pulls case counts per country and re-sorts the table (rows
# 4..204) descending by column B (Casos totales). Because the sheet is
# sorted, a changed total can shuffle which country occupies which row,
# which is why the diff shows shared-string reordering even though most
# rows keep the same numbers. Re-creating that sort here (by just writing
# the correct country name + values into each affected row) reproduces the
# same effect without needing to touch the shared-string table directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (row 1, A1).
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 00:12"

# Rows whose country and/or numbers change as a result of the data refresh
# + re-sort. Each entry: row, País, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos críticos, Muertes hoy, Muertes.
$rowsData = @(
    @(6, 'Estados Unidos', 65652, 10796, 394, 64327, 1411, 151, 931),
    @(20, 'Australia', 2676, 359, 118, 2547, 11, 3, 11),
    @(21, 'Suecia', 2526, 227, 16, 2448, 158, 22, 62),
    @(22, 'Brasil', 2433, 186, 2, 2374, 18, 11, 57),
    @(23, 'Turquia', 2433, 561, 26, 2348, 136, 15, 59),
    @(116, 'Consejo Danes para los Refugiados', 48, 3, 0, 46, 0, 0, 2),
    @(117, 'Mauricio', 48, 6, 0, 46, 1, 0, 2),
    @(123, 'Honduras', 36, 6, 0, 36, 0, 0, 0),
    @(124, 'Mayotte', 36, 0, 0, 36, 0, 0, 0),
    @(141, 'Uganda', 14, 5, 0, 14, 0, 0, 0),
    @(142, 'Nueva Caledonia', 14, 4, 0, 14, 0, 0, 0),
    @(145, 'Etiopia', 12, 0, 0, 12, 0, 0, 0),
    @(146, 'Tanzania', 12, 0, 0, 12, 0, 0, 0),
    @(147, 'Republica de Yibuti', 11, 8, 0, 11, 0, 0, 0),
    @(148, 'San Martin (Parte Francesa)', 11, 3, 0, 11, 0, 0, 0),
    @(152, 'Haiti', 8, 1, 0, 8, 0, 0, 0),
    @(153, 'Surinam', 8, 1, 0, 8, 0, 0, 0),
    @(164, 'Siria', 5, 4, 0, 5, 0, 0, 0),
    @(165, 'Mozambique', 5, 2, 0, 5, 0, 0, 0),
    @(169, 'Congo', 4, 0, 0, 4, 0, 0, 0),
    @(171, 'Eritrea', 4, 3, 0, 4, 0, 0, 0),
    @(175, 'Santa Lucia', 3, 0, 0, 3, 0, 0, 0),
    @(176, 'Birmania', 3, 0, 0, 3, 0, 0, 0),
    @(179, 'Angola', 3, 0, 0, 3, 0, 0, 0),
    @(180, 'Republica de Africa Central', 3, 0, 0, 3, 0, 0, 0),
    @(181, 'Republica del Chad', 3, 0, 0, 3, 0, 0, 0),
    @(182, 'San Bartolome', 3, 0, 0, 3, 0, 0, 0),
    @(183, 'Liberia', 3, 0, 0, 3, 0, 0, 0),
    @(184, 'Laos', 3, 1, 0, 3, 0, 0, 0),
    @(189, 'Islas Virgenes Britanicas', 2, 2, 0, 2, 0, 0, 0),
    @(191, 'Nicaragua', 2, 0, 0, 2, 0, 0, 0),
    @(192, 'Belice', 2, 1, 0, 2, 0, 0, 0),
    @(193, 'Butan', 2, 0, 0, 2, 0, 0, 0),
    @(194, 'Mauritania', 2, 0, 0, 2, 0, 0, 0),
    @(196, 'Guinea-Bisau', 2, 2, 0, 2, 0, 0, 0),
    @(197, 'Papua Nueva Guinea', 1, 0, 0, 1, 0, 0, 0),
    @(198, 'Libia', 1, 0, 0, 1, 0, 0, 0),
    @(199, 'Somalia', 1, 0, 0, 1, 0, 0, 0),
    @(200, 'Montserrat', 1, 0, 0, 1, 0, 0, 0),
    @(201, 'San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0),
    @(202, 'Granada', 1, 0, 0, 1, 0, 0, 0),
    @(203, 'Timor Oriental', 1, 0, 0, 1, 0, 0, 0),
    @(204, 'Islas Turcas y Caicos', 1, 0, 0, 1, 0, 0, 0)
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
    $ws.Cells.Item($r, 6).Value = $entry[6]
    $ws.Cells.Item($r, 7).Value = $entry[7]
    $ws.Cells.Item($r, 8).Value = $entry[8]
}
